$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.434.30'
$ws.Range('E2').Value = '  -1.27%  '

$ws.Range('D3').Value = '2.514.77'
$ws.Range('E3').Value = '  -0.09%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').Value = '572.73'
$ws.Range('E5').Value = '  -0.28%  '

$ws.Range('D6').Value = '166.73'
$ws.Range('E6').Value = '  -1.36%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('E8').Value = '  +2.14%  '

$ws.Range('D9').Value = '2.513.83'
$ws.Range('E9').Value = '  -0.11%  '

$ws.Range('E10').Value = '  -1.44%  '

$ws.Range('E11').Value = '  -0.58%  '

$ws.Range('E12').Value = '  +4.26%  '

$ws.Range('D13').Value = '4.93'
$ws.Range('E13').Value = '  +3.33%  '

$ws.Range('D14').Value = '2.976.97'
$ws.Range('E14').Value = '  -0.01%  '

$ws.Range('D15').Value = '69.332.30'
$ws.Range('E15').Value = '  -1.27%  '

$ws.Range('D16').Value = '0.0000175'
$ws.Range('E16').Value = '  -2.22%  '

$ws.Range('D17').Value = '24.89'
$ws.Range('E17').Value = '  +0.33%  '

$ws.Range('D18').Value = '2.513.57'
$ws.Range('E18').Value = '  +0.00%  '

$ws.Range('D19').Value = '11.35'
$ws.Range('E19').Value = '  -1.27%  '

$ws.Range('D20').Value = '7.71'
$ws.Range('E20').Value = '  +2.39%  '

$ws.Range('D21').Value = '349.36'
$ws.Range('E21').Value = '  -1.67%  '

$ws.Range('E22').Value = '  +0.47%  '

$ws.Range('D23').Value = '1.97'
$ws.Range('E23').Value = '  +1.28%  '

$ws.Range('E24').Value = '  +0.05%  '

$ws.Range('D25').Value = '70.01'

$ws.Range('E26').Value = '  -1.20%  '

$ws.Range('D27').Value = '8.92'
$ws.Range('E27').Value = '  -2.80%  '

$ws.Range('D28').Value = '2.648.28'
$ws.Range('E28').Value = '  -0.01%  '

$ws.Range('E29').Value = '  +0.03%  '

$ws.Range('D30').Value = '0.0₃0893'
$ws.Range('E30').Value = '  -1.36%  '

$ws.Range('D31').Value = '7.88'
$ws.Range('E31').Value = '  +0.85%  '

$ws.Range('D32').Value = '463.96'
$ws.Range('E32').Value = '  -2.83%  '

$ws.Range('E33').Value = '  -1.17%  '

$ws.Range('E34').Value = '  -0.90%  '

$ws.Range('D35').Value = '0.999'

$ws.Range('E36').Value = '  +1.05%  '

$ws.Range('D37').Value = '157.76'
$ws.Range('E37').Value = '  +0.47%  '

$ws.Range('D38').Value = '19.02'
$ws.Range('E38').Value = '  +1.11%  '

$ws.Range('D39').Value = '18.53'
$ws.Range('E39').Value = '  +0.12%  '

$ws.Range('E40').Value = '  -0.03%  '

$ws.Range('E41').Value = '  +1.49%  '

$ws.Range('D42').Value = '0.319'
$ws.Range('E42').Value = '  +0.63%  '

$ws.Range('E43').Value = '  -2.14%  '

$ws.Range('D44').Value = '38.32'
$ws.Range('E44').Value = '  +0.19%  '

$ws.Range('D45').Value = '2.27'
$ws.Range('E45').Value = '  -5.21%  '

$ws.Range('D46').Value = '1.13'
$ws.Range('E46').Value = '  -12.75%  '

$ws.Range('D47').Value = '141.59'
$ws.Range('E47').Value = '  -0.51%  '

$ws.Range('E48').Value = '  +1.13%  '

$ws.Range('D49').Value = '3.48'
$ws.Range('E49').Value = '  -0.77%  '

$ws.Range('E50').Value = '  -0.62%  '

